$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 4800
$ws.Range("I38").Value = 3400
$ws.Range("J38").Value = 9000
$ws.Range("K38").Value = 10200
$ws.Range("L38").Value = 27000
$ws.Range("M38").Value = -9828
$ws.Range("N38").Value = -27744
$ws.Range("H41").Value = 388.6
$ws.Range("I41").Value = 18.333334
$ws.Range("K41").Value = 18.333334
$ws.Range("M41").Value = 421.666666
$ws.Range("H43").Value = 2499
$ws.Range("J43").Value = 2647.5715
$ws.Range("L43").Value = 2647.5715
$ws.Range("N43").Value = -2785.5715
$ws.Range("H49").Value = 435
$ws.Range("I49").Value = 435
$ws.Range("K49").Value = 1305
$ws.Range("M49").Value = -1169
$ws.Range("H55").Value = 492.53845
$ws.Range("I55").Value = 307.8
$ws.Range("J55").Value = 608
$ws.Range("K55").Value = 307.8
$ws.Range("L55").Value = 608
$ws.Range("M55").Value = -93.80000000000001
$ws.Range("N55").Value = -1036
$ws.Range("H64").Value = 5315.857
$ws.Range("I64").Value = 4702.923
$ws.Range("J64").Value = 6311.875
$ws.Range("K64").Value = 4702.923
$ws.Range("L64").Value = 6311.875
$ws.Range("M64").Value = -4454.923
$ws.Range("N64").Value = -6807.875
$ws.Range("H67").Value = 5315.857
$ws.Range("I67").Value = 4702.923
$ws.Range("J67").Value = 6311.875
$ws.Range("K67").Value = 4702.923
$ws.Range("L67").Value = 6311.875
$ws.Range("M67").Value = -3844.923
$ws.Range("N67").Value = -8027.875
$ws.Range("H76").Value = 10762.895
$ws.Range("I76").Value = 10617.412
$ws.Range("K76").Value = 10617.412
$ws.Range("M76").Value = -10302.412
$ws.Range("H79").Value = 10762.895
$ws.Range("I79").Value = 10617.412
$ws.Range("K79").Value = 10617.412
$ws.Range("M79").Value = -9525.412
$ws.Range("H80").Value = 1784.25
$ws.Range("I80").Value = 1745.8334
$ws.Range("J80").Value = 1899.5
$ws.Range("K80").Value = 5237.5002
$ws.Range("L80").Value = 5698.5
$ws.Range("M80").Value = -4239.5002
$ws.Range("N80").Value = -7694.5
$ws.Range("H83").Value = 1784.25
$ws.Range("I83").Value = 1745.8334
$ws.Range("J83").Value = 1899.5
$ws.Range("K83").Value = 15712.5006
$ws.Range("L83").Value = 17095.5
$ws.Range("M83").Value = -10720.5006
$ws.Range("N83").Value = -27079.5
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384
$ws.Range("H92").Value = 4776.5835
$ws.Range("J92").Value = 7975.8
$ws.Range("L92").Value = 7975.8
$ws.Range("N92").Value = -10471.8
$ws.Range("H96").Value = 463.7143
$ws.Range("I96").Value = 577.4
$ws.Range("K96").Value = 1732.2
$ws.Range("M96").Value = -359.1999999999998
$ws.Range("H97").Value = 2799.6667
$ws.Range("J97").Value = 2799.6667
$ws.Range("L97").Value = 8399.000100000001
$ws.Range("N97").Value = -9391.000100000001
$ws.Range("H100").Value = 1074.25
$ws.Range("I100").Value = 1074.25
$ws.Range("K100").Value = 1074.25
$ws.Range("M100").Value = -533.25
$ws.Range("H106").Value = 7064.3335
$ws.Range("I106").Value = 7064.3335
$ws.Range("K106").Value = 7064.3335
$ws.Range("M106").Value = -6433.3335
$ws.Range("H111").Value = 3325
$ws.Range("I111").Value = 3124.25
$ws.Range("J111").Value = 3485.6
$ws.Range("K111").Value = 9372.75
$ws.Range("L111").Value = 10456.8
$ws.Range("M111").Value = -6305.75
$ws.Range("N111").Value = -16590.8
$ws.Range("H112").Value = 1808.742
$ws.Range("I112").Value = 1300
$ws.Range("J112").Value = 1863.25
$ws.Range("K112").Value = 3900
$ws.Range("L112").Value = 5589.75
$ws.Range("M112").Value = -2792
$ws.Range("N112").Value = -7805.75
$ws.Range("H113").Value = 6248.9
$ws.Range("I113").Value = 6937.375
$ws.Range("J113").Value = 3495
$ws.Range("K113").Value = 6937.375
$ws.Range("L113").Value = 3495
$ws.Range("M113").Value = -3683.375
$ws.Range("N113").Value = -10003
$ws.Range("H137").Value = 7275.3335
$ws.Range("I137").Value = 4629
$ws.Range("K137").Value = 13887
$ws.Range("M137").Value = -11337
$ws.Range("H138").Value = 1975.2858
$ws.Range("I138").Value = 2018
$ws.Range("J138").Value = 1955.0526
$ws.Range("K138").Value = 6054
$ws.Range("L138").Value = 5865.1578
$ws.Range("M138").Value = -914
$ws.Range("N138").Value = -16145.1578
$ws.Range("H141").Value = 5055
$ws.Range("I141").Value = 2220
$ws.Range("K141").Value = 6660
$ws.Range("M141").Value = -1480

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 1141.5
$ws.Range("I16").Value = 1141.5
$ws.Range("K16").Value = 1141.5
$ws.Range("M16").Value = -854.5
$ws.Range("H32").Value = 10746.3125
$ws.Range("I32").Value = 6462.3
$ws.Range("J32").Value = 17886.334
$ws.Range("K32").Value = 6462.3
$ws.Range("L32").Value = 17886.334
$ws.Range("M32").Value = -6175.3
$ws.Range("N32").Value = -18460.334
$ws.Range("H45").Value = 1990
$ws.Range("I45").Value = 1990
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1990
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1613
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 4859.5
$ws.Range("I61").Value = 4656
$ws.Range("J61").Value = 4900.2
$ws.Range("K61").Value = 4656
$ws.Range("L61").Value = 4900.2
$ws.Range("M61").Value = -4444
$ws.Range("N61").Value = -5324.2
$ws.Range("H74").Value = 1678.6
$ws.Range("I74").Value = 1686
$ws.Range("K74").Value = 1686
$ws.Range("M74").Value = -812
$ws.Range("H77").Value = 1678.6
$ws.Range("I77").Value = 1686
$ws.Range("K77").Value = 8430
$ws.Range("M77").Value = -4062
$ws.Range("H80").Value = 54666.668
$ws.Range("J80").Value = 54666.668
$ws.Range("L80").Value = 54666.668
$ws.Range("N80").Value = -56662.668
$ws.Range("H83").Value = 54666.668
$ws.Range("J83").Value = 54666.668
$ws.Range("L83").Value = 164000.004
$ws.Range("N83").Value = -173984.004
$ws.Range("H102").Value = 1822.8948
$ws.Range("I102").Value = 1656.5294
$ws.Range("K102").Value = 1656.5294
$ws.Range("M102").Value = -34.5293999999999
$ws.Range("H110").Value = 3550.1
$ws.Range("I110").Value = 2981.8
$ws.Range("K110").Value = 2981.8
$ws.Range("M110").Value = -936.8000000000002
$ws.Range("H118").Value = 60000
$ws.Range("J118").Value = 60000
$ws.Range("L118").Value = 60000
$ws.Range("N118").Value = -63314
$ws.Range("H122").Value = 2599.4827
$ws.Range("I122").Value = 2518.4285
$ws.Range("J122").Value = 2812.25
$ws.Range("K122").Value = 7555.2855
$ws.Range("L122").Value = 8436.75
$ws.Range("M122").Value = -5105.2855
$ws.Range("N122").Value = -13336.75
$ws.Range("H132").Value = 5401.8
$ws.Range("J132").Value = 3999
$ws.Range("L132").Value = 11997
$ws.Range("N132").Value = -17057
$ws.Range("H134").Value = 19998
$ws.Range("J134").Value = 19998
$ws.Range("L134").Value = 19998
$ws.Range("N134").Value = -30138
$ws.Range("H136").Value = 4859.5
$ws.Range("I136").Value = 4656
$ws.Range("J136").Value = 4900.2
$ws.Range("K136").Value = 13968
$ws.Range("L136").Value = 14700.6
$ws.Range("M136").Value = -11418
$ws.Range("N136").Value = -19800.6
$ws.Range("H140").Value = 67500
$ws.Range("J140").Value = 67500
$ws.Range("L140").Value = 67500
$ws.Range("N140").Value = -77860

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7500
$ws.Range("I20").Value = 7500
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 7500
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -7253
$ws.Range("N20").ClearContents()
$ws.Range("H82").Value = 44944.223
$ws.Range("I82").Value = 15166.333
$ws.Range("J82").Value = 59833.168
$ws.Range("K82").Value = 15166.333
$ws.Range("L82").Value = 59833.168
$ws.Range("M82").Value = -14783.333
$ws.Range("N82").Value = -60599.168
$ws.Range("H85").Value = 44944.223
$ws.Range("I85").Value = 15166.333
$ws.Range("J85").Value = 59833.168
$ws.Range("K85").Value = 15166.333
$ws.Range("L85").Value = 59833.168
$ws.Range("M85").Value = -13840.333
$ws.Range("N85").Value = -62485.168
$ws.Range("H86").Value = 2975.2222
$ws.Range("I86").Value = 2722.125
$ws.Range("K86").Value = 2722.125
$ws.Range("M86").Value = -1599.125
$ws.Range("H88").Value = 30624.25
$ws.Range("J88").Value = 32142
$ws.Range("L88").Value = 32142
$ws.Range("N88").Value = -32954
$ws.Range("H89").Value = 2975.2222
$ws.Range("I89").Value = 2722.125
$ws.Range("K89").Value = 13610.625
$ws.Range("M89").Value = -7994.625
$ws.Range("H91").Value = 30624.25
$ws.Range("J91").Value = 32142
$ws.Range("L91").Value = 32142
$ws.Range("N91").Value = -34950
$ws.Range("H94").Value = 2356.0417
$ws.Range("I94").Value = 1291.2354
$ws.Range("K94").Value = 1291.2354
$ws.Range("M94").Value = -840.2354
$ws.Range("H99").Value = 4283.923
$ws.Range("I99").Value = 4156.143
$ws.Range("J99").Value = 4433
$ws.Range("K99").Value = 4156.143
$ws.Range("L99").Value = 4433
$ws.Range("M99").Value = -2658.143
$ws.Range("N99").Value = -7429
$ws.Range("H105").Value = 4610.1304
$ws.Range("I105").Value = 3574.2942
$ws.Range("K105").Value = 3574.2942
$ws.Range("M105").Value = -1827.2942
$ws.Range("H107").Value = 3000
$ws.Range("I107").Value = 2000
$ws.Range("K107").Value = 2000
$ws.Range("M107").Value = -80
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6852.2085
$ws.Range("I31").Value = 3111
$ws.Range("J31").Value = 8099.278
$ws.Range("K31").Value = 3111
$ws.Range("L31").Value = 8099.278
$ws.Range("M31").Value = -2816
$ws.Range("N31").Value = -8689.278
$ws.Range("H34").Value = 6852.2085
$ws.Range("I34").Value = 3111
$ws.Range("J34").Value = 8099.278
$ws.Range("K34").Value = 3111
$ws.Range("L34").Value = 8099.278
$ws.Range("M34").Value = -2909
$ws.Range("N34").Value = -8503.278
$ws.Range("H55").Value = 56000
$ws.Range("J55").Value = 56000
$ws.Range("L55").Value = 56000
$ws.Range("N55").Value = -56630
$ws.Range("H58").Value = 7347.5
$ws.Range("I58").Value = 7497.222
$ws.Range("K58").Value = 7497.222
$ws.Range("M58").Value = -7294.222
$ws.Range("H86").Value = 14998.5
$ws.Range("J86").Value = 14998
$ws.Range("L86").Value = 14998
$ws.Range("N86").Value = -17244
$ws.Range("H89").Value = 14998.5
$ws.Range("J89").Value = 14998
$ws.Range("L89").Value = 74990
$ws.Range("N89").Value = -86222
$ws.Range("H132").Value = 997.25
$ws.Range("I132").Value = 997.25
$ws.Range("K132").Value = 2991.75
$ws.Range("M132").Value = -461.75
$ws.Range("H134").Value = 8370.454
$ws.Range("I134").Value = 8155.6313
$ws.Range("K134").Value = 24466.8939
$ws.Range("M134").Value = -21931.8939
$ws.Range("H136").Value = 7347.5
$ws.Range("I136").Value = 7497.222
$ws.Range("K136").Value = 22491.666
$ws.Range("M136").Value = -19941.666
$ws.Range("H140").Value = 99996.5
$ws.Range("J140").Value = 99996.5
$ws.Range("L140").Value = 99996.5
$ws.Range("N140").Value = -110356.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 250708.75
$ws.Range("J32").Value = 945
$ws.Range("L32").Value = 2835
$ws.Range("N32").Value = -3401
$ws.Range("H33").Value = 2200
$ws.Range("I33").Value = 744.4
$ws.Range("K33").Value = 4466.4
$ws.Range("M33").Value = -4183.4
$ws.Range("H39").Value = 2866.3333
$ws.Range("J39").Value = 2866.3333
$ws.Range("L39").Value = 8598.999899999999
$ws.Range("N39").Value = -9186.999899999999
$ws.Range("H80").Value = 15527.286
$ws.Range("I80").Value = 34166.332
$ws.Range("J80").Value = 1548
$ws.Range("K80").Value = 102498.996
$ws.Range("L80").Value = 4644
$ws.Range("M80").Value = -101562.996
$ws.Range("N80").Value = -6516
$ws.Range("H83").Value = 15527.286
$ws.Range("I83").Value = 34166.332
$ws.Range("J83").Value = 1548
$ws.Range("K83").Value = 307496.988
$ws.Range("L83").Value = 13932
$ws.Range("M83").Value = -302816.988
$ws.Range("N83").Value = -23292
$ws.Range("H107").Value = 1096.5
$ws.Range("I107").Value = 194.33333
$ws.Range("K107").Value = 582.99999
$ws.Range("M107").Value = 1337.00001
$ws.Range("H111").Value = 1300
$ws.Range("I111").Value = 1300
$ws.Range("K111").Value = 3900
$ws.Range("M111").Value = -833
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("H137").Value = 5444.933
$ws.Range("I137").Value = 7347.5557
$ws.Range("K137").Value = 22042.6671
$ws.Range("M137").Value = -16942.6671
$ws.Range("H138").Value = 2321
$ws.Range("I138").Value = 2321
$ws.Range("K138").Value = 6963
$ws.Range("M138").Value = -1823

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7553.909
$ws.Range("I70").Value = 4448.5
$ws.Range("K70").Value = 4448.5
$ws.Range("M70").Value = -4178.5
$ws.Range("H73").Value = 7553.909
$ws.Range("I73").Value = 4448.5
$ws.Range("K73").Value = 4448.5
$ws.Range("M73").Value = -3512.5
$ws.Range("H80").Value = 9918.531999999999
$ws.Range("I80").Value = 6823.2
$ws.Range("J80").Value = 12211.37
$ws.Range("K80").Value = 6823.2
$ws.Range("L80").Value = 12211.37
$ws.Range("M80").Value = -5825.2
$ws.Range("N80").Value = -14207.37
$ws.Range("H83").Value = 9918.531999999999
$ws.Range("I83").Value = 6823.2
$ws.Range("J83").Value = 12211.37
$ws.Range("K83").Value = 34116
$ws.Range("L83").Value = 61056.85000000001
$ws.Range("M83").Value = -29124
$ws.Range("N83").Value = -71040.85000000001
$ws.Range("H97").Value = 831.2273
$ws.Range("I97").Value = 843.2778
$ws.Range("J97").Value = 777
$ws.Range("K97").Value = 843.2778
$ws.Range("L97").Value = 777
$ws.Range("M97").Value = -347.2778
$ws.Range("N97").Value = -1769
$ws.Range("H113").Value = 1835.8948
$ws.Range("I113").Value = 1711.5625
$ws.Range("K113").Value = 1711.5625
$ws.Range("M113").Value = 458.4375
$ws.Range("H122").Value = 79418.30499999999
$ws.Range("I122").Value = 168583.17
$ws.Range("K122").Value = 505749.51
$ws.Range("M122").Value = -503299.51
$ws.Range("H126").Value = 4497.625
$ws.Range("I126").Value = 4200
$ws.Range("K126").Value = 12600
$ws.Range("M126").Value = -10130
$ws.Range("H132").Value = 4399.2
$ws.Range("I132").Value = 4399.2
$ws.Range("K132").Value = 13197.6
$ws.Range("M132").Value = -10667.6
$ws.Range("H134").Value = 49500
$ws.Range("J134").Value = 49500
$ws.Range("L134").Value = 148500
$ws.Range("N134").Value = -153570

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6381.5713
$ws.Range("I7").Value = 7261.25
$ws.Range("J7").Value = 5840.231
$ws.Range("K7").Value = 7261.25
$ws.Range("L7").Value = 5840.231
$ws.Range("M7").Value = -7149.25
$ws.Range("N7").Value = -6064.231
$ws.Range("H40").Value = 4121.9443
$ws.Range("I40").Value = 3541
$ws.Range("J40").Value = 4848.125
$ws.Range("K40").Value = 3541
$ws.Range("L40").Value = 4848.125
$ws.Range("M40").Value = -3405
$ws.Range("N40").Value = -5120.125
$ws.Range("H68").Value = 1884.5385
$ws.Range("I68").Value = 2045.3636
$ws.Range("J68").Value = 1000
$ws.Range("K68").Value = 2045.3636
$ws.Range("L68").Value = 1000
$ws.Range("M68").Value = -1296.3636
$ws.Range("N68").Value = -2498
$ws.Range("H71").Value = 1884.5385
$ws.Range("I71").Value = 2045.3636
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 10226.818
$ws.Range("L71").Value = 5000
$ws.Range("M71").Value = -6482.817999999999
$ws.Range("N71").Value = -12488
$ws.Range("H93").Value = 2493.5417
$ws.Range("I93").Value = 2372
$ws.Range("J93").Value = 3830.5
$ws.Range("K93").Value = 2372
$ws.Range("L93").Value = 3830.5
$ws.Range("M93").Value = -1124
$ws.Range("N93").Value = -6326.5
$ws.Range("H126").Value = 6381.5713
$ws.Range("I126").Value = 7261.25
$ws.Range("J126").Value = 5840.231
$ws.Range("K126").Value = 21783.75
$ws.Range("L126").Value = 17520.693
$ws.Range("M126").Value = -19313.75
$ws.Range("N126").Value = -22460.693
$ws.Range("H132").Value = 3218.7
$ws.Range("I132").Value = 3240.6667
$ws.Range("K132").Value = 9722.000100000001
$ws.Range("M132").Value = -7192.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3263.9375
$ws.Range("I81").Value = 3214.8667
$ws.Range("J81").Value = 4000
$ws.Range("K81").Value = 6429.7334
$ws.Range("L81").Value = 8000
$ws.Range("M81").Value = -5368.7334
$ws.Range("N81").Value = -10122
$ws.Range("H84").Value = 3263.9375
$ws.Range("I84").Value = 3214.8667
$ws.Range("J84").Value = 4000
$ws.Range("K84").Value = 32148.667
$ws.Range("L84").Value = 40000
$ws.Range("M84").Value = -26844.667
$ws.Range("N84").Value = -50608
$ws.Range("H96").Value = 3527.3635
$ws.Range("I96").Value = 3580
$ws.Range("J96").Value = 3387
$ws.Range("K96").Value = 3580
$ws.Range("L96").Value = 3387
$ws.Range("M96").Value = -2207
$ws.Range("N96").Value = -6133
$ws.Range("H100").Value = 725.6667
$ws.Range("I100").Value = 725.6667
$ws.Range("K100").Value = 1451.3334
$ws.Range("M100").Value = -910.3334
$ws.Range("H107").Value = 2402.2
$ws.Range("J107").Value = 2752.5
$ws.Range("L107").Value = 8257.5
$ws.Range("N107").Value = -12097.5
$ws.Range("H113").Value = 718.75
$ws.Range("I113").Value = 688.3333
$ws.Range("J113").Value = 810
$ws.Range("K113").Value = 2064.9999
$ws.Range("L113").Value = 2430
$ws.Range("M113").Value = 105.0001000000002
$ws.Range("N113").Value = -6770
$ws.Range("H122").Value = 5439.727
$ws.Range("I122").Value = 4935.9287
$ws.Range("J122").Value = 8261
$ws.Range("K122").Value = 14807.7861
$ws.Range("L122").Value = 24783
$ws.Range("M122").Value = -12357.7861
$ws.Range("N122").Value = -29683
$ws.Range("H126").Value = 1503.8695
$ws.Range("I126").Value = 1435.9546
$ws.Range("K126").Value = 4307.8638
$ws.Range("M126").Value = -1837.8638
$ws.Range("H132").Value = 2404.077
$ws.Range("I132").Value = 2404.077
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7212.231000000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4682.231000000001
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 4271.5
$ws.Range("I136").Value = 3962
$ws.Range("K136").Value = 11886
$ws.Range("M136").Value = -9336
